$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update packaging material names for rows 15 and 16 (previously duplicated as "Non woven bag 1")
$ws.Range("A15").Value = "Non woven bag 2"
$ws.Range("A16").Value = "Non woven bag 3"

# Move the active selection/cursor to K10, matching the final saved cursor position
$ws.Range("K10").Select()
